$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V between paired rows (home/away reordering corrections) ---
$tmp = $ws.Range("F12:V12").Value2
$ws.Range("F12:V12").Value2 = $ws.Range("F13:V13").Value2
$ws.Range("F13:V13").Value2 = $tmp
$tmp = $ws.Range("F15:V15").Value2
$ws.Range("F15:V15").Value2 = $ws.Range("F16:V16").Value2
$ws.Range("F16:V16").Value2 = $tmp
$tmp = $ws.Range("F23:V23").Value2
$ws.Range("F23:V23").Value2 = $ws.Range("F24:V24").Value2
$ws.Range("F24:V24").Value2 = $tmp
$tmp = $ws.Range("F26:V26").Value2
$ws.Range("F26:V26").Value2 = $ws.Range("F27:V27").Value2
$ws.Range("F27:V27").Value2 = $tmp
$tmp = $ws.Range("F46:V46").Value2
$ws.Range("F46:V46").Value2 = $ws.Range("F47:V47").Value2
$ws.Range("F47:V47").Value2 = $tmp
$tmp = $ws.Range("F49:V49").Value2
$ws.Range("F49:V49").Value2 = $ws.Range("F50:V50").Value2
$ws.Range("F50:V50").Value2 = $tmp
$tmp = $ws.Range("F51:V51").Value2
$ws.Range("F51:V51").Value2 = $ws.Range("F52:V52").Value2
$ws.Range("F52:V52").Value2 = $tmp
$tmp = $ws.Range("F53:V53").Value2
$ws.Range("F53:V53").Value2 = $ws.Range("F54:V54").Value2
$ws.Range("F54:V54").Value2 = $tmp
$tmp = $ws.Range("F55:V55").Value2
$ws.Range("F55:V55").Value2 = $ws.Range("F56:V56").Value2
$ws.Range("F56:V56").Value2 = $tmp
$tmp = $ws.Range("F63:V63").Value2
$ws.Range("F63:V63").Value2 = $ws.Range("F64:V64").Value2
$ws.Range("F64:V64").Value2 = $tmp
$tmp = $ws.Range("F75:V75").Value2
$ws.Range("F75:V75").Value2 = $ws.Range("F76:V76").Value2
$ws.Range("F76:V76").Value2 = $tmp
$tmp = $ws.Range("F87:V87").Value2
$ws.Range("F87:V87").Value2 = $ws.Range("F88:V88").Value2
$ws.Range("F88:V88").Value2 = $tmp

# --- Append 10 new match rows (91-100) ---
# Row 91
$ws.Range("A90:V90").Copy($ws.Range("A91"))
$ws.Range("B91").Value2 = "turkey"
$ws.Range("C91").Value2 = "super-lig"
$ws.Range("D91").Value2 = "2023-2024"
$ws.Range("A91").Value2 = 90
$ws.Range("E91").Value2 = 45226.79166666666
$ws.Range("F91").Value2 = "Hatayspor"
$ws.Range("G91").Value2 = 1
$ws.Range("H91").Value2 = "Kayserispor"
$ws.Range("I91").Value2 = 2
$ws.Range("J91").Value2 = 2.26
$ws.Range("K91").Value2 = "22/10/2023 20:15"
$ws.Range("L91").Value2 = 2.49
$ws.Range("M91").Value2 = "27/10/2023 18:56"
$ws.Range("N91").Value2 = 3.55
$ws.Range("O91").Value2 = "22/10/2023 20:15"
$ws.Range("P91").Value2 = 3.42
$ws.Range("Q91").Value2 = "27/10/2023 18:56"
$ws.Range("R91").Value2 = 3.26
$ws.Range("S91").Value2 = "22/10/2023 20:15"
$ws.Range("T91").Value2 = 3.01
$ws.Range("U91").Value2 = "27/10/2023 18:56"
$ws.Range("V91").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/hatayspor-kayserispor/069EREiq/"

# Row 92
$ws.Range("A91:V91").Copy($ws.Range("A92"))
$ws.Range("B92").Value2 = "turkey"
$ws.Range("C92").Value2 = "super-lig"
$ws.Range("D92").Value2 = "2023-2024"
$ws.Range("A92").Value2 = 91
$ws.Range("E92").Value2 = 45226.79166666666
$ws.Range("F92").Value2 = "Kasimpasa"
$ws.Range("G92").Value2 = 3
$ws.Range("H92").Value2 = "Istanbulspor AS"
$ws.Range("I92").Value2 = 1
$ws.Range("J92").Value2 = 1.86
$ws.Range("K92").Value2 = "22/10/2023 15:12"
$ws.Range("L92").Value2 = 1.96
$ws.Range("M92").Value2 = "27/10/2023 18:58"
$ws.Range("N92").Value2 = 3.96
$ws.Range("O92").Value2 = "22/10/2023 15:12"
$ws.Range("P92").Value2 = 3.71
$ws.Range("Q92").Value2 = "27/10/2023 18:59"
$ws.Range("R92").Value2 = 4.09
$ws.Range("S92").Value2 = "22/10/2023 15:12"
$ws.Range("T92").Value2 = 3.99
$ws.Range("U92").Value2 = "27/10/2023 18:58"
$ws.Range("V92").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-istanbulspor-as/YJ7MPhMe/"

# Row 93
$ws.Range("A92:V92").Copy($ws.Range("A93"))
$ws.Range("B93").Value2 = "turkey"
$ws.Range("C93").Value2 = "super-lig"
$ws.Range("D93").Value2 = "2023-2024"
$ws.Range("A93").Value2 = 92
$ws.Range("E93").Value2 = 45227.625
$ws.Range("F93").Value2 = "Alanyaspor"
$ws.Range("G93").Value2 = 1
$ws.Range("H93").Value2 = "Sivasspor"
$ws.Range("I93").Value2 = 2
$ws.Range("J93").Value2 = 2.1
$ws.Range("K93").Value2 = "23/10/2023 19:12"
$ws.Range("L93").Value2 = 2.26
$ws.Range("M93").Value2 = "28/10/2023 14:57"
$ws.Range("N93").Value2 = 3.64
$ws.Range("O93").Value2 = "23/10/2023 19:12"
$ws.Range("P93").Value2 = 3.34
$ws.Range("Q93").Value2 = "28/10/2023 14:58"
$ws.Range("R93").Value2 = 3.55
$ws.Range("S93").Value2 = "23/10/2023 19:12"
$ws.Range("T93").Value2 = 3.49
$ws.Range("U93").Value2 = "28/10/2023 14:56"
$ws.Range("V93").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/alanyaspor-sivasspor/W4QPOCy2/"

# Row 94
$ws.Range("A93:V93").Copy($ws.Range("A94"))
$ws.Range("B94").Value2 = "turkey"
$ws.Range("C94").Value2 = "super-lig"
$ws.Range("D94").Value2 = "2023-2024"
$ws.Range("A94").Value2 = 93
$ws.Range("E94").Value2 = 45227.75
$ws.Range("F94").Value2 = "Adana Demirspor"
$ws.Range("G94").Value2 = 3
$ws.Range("H94").Value2 = "Konyaspor"
$ws.Range("I94").Value2 = 0
$ws.Range("J94").Value2 = 1.73
$ws.Range("K94").Value2 = "23/10/2023 19:12"
$ws.Range("L94").Value2 = 1.83
$ws.Range("M94").Value2 = "28/10/2023 17:58"
$ws.Range("N94").Value2 = 4.22
$ws.Range("O94").Value2 = "23/10/2023 19:12"
$ws.Range("P94").Value2 = 3.99
$ws.Range("Q94").Value2 = "28/10/2023 17:59"
$ws.Range("R94").Value2 = 4.54
$ws.Range("S94").Value2 = "23/10/2023 19:12"
$ws.Range("T94").Value2 = 4.27
$ws.Range("U94").Value2 = "28/10/2023 17:59"
$ws.Range("V94").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/adanademirspor-konyaspor/Y948mdE2/"

# Row 95
$ws.Range("A94:V94").Copy($ws.Range("A95"))
$ws.Range("B95").Value2 = "turkey"
$ws.Range("C95").Value2 = "super-lig"
$ws.Range("D95").Value2 = "2023-2024"
$ws.Range("A95").Value2 = 94
$ws.Range("E95").Value2 = 45227.75
$ws.Range("F95").Value2 = "Rizespor"
$ws.Range("G95").Value2 = 0
$ws.Range("H95").Value2 = "Galatasaray"
$ws.Range("I95").Value2 = 1
$ws.Range("J95").Value2 = 5.1
$ws.Range("K95").Value2 = "23/10/2023 05:42"
$ws.Range("L95").Value2 = 7.86
$ws.Range("M95").Value2 = "28/10/2023 17:37"
$ws.Range("N95").Value2 = 4.48
$ws.Range("O95").Value2 = "23/10/2023 05:42"
$ws.Range("P95").Value2 = 5.37
$ws.Range("Q95").Value2 = "28/10/2023 17:37"
$ws.Range("R95").Value2 = 1.6
$ws.Range("S95").Value2 = "23/10/2023 05:42"
$ws.Range("T95").Value2 = 1.39
$ws.Range("U95").Value2 = "28/10/2023 17:37"
$ws.Range("V95").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/rizespor-galatasaray/vBmiuYLr/"

# Row 96
$ws.Range("A95:V95").Copy($ws.Range("A96"))
$ws.Range("B96").Value2 = "turkey"
$ws.Range("C96").Value2 = "super-lig"
$ws.Range("D96").Value2 = "2023-2024"
$ws.Range("A96").Value2 = 95
$ws.Range("E96").Value2 = 45228.58333333334
$ws.Range("F96").Value2 = "Karagumruk"
$ws.Range("G96").Value2 = 0
$ws.Range("H96").Value2 = "Trabzonspor"
$ws.Range("I96").Value2 = 0
$ws.Range("J96").Value2 = 2.63
$ws.Range("K96").Value2 = "23/10/2023 19:12"
$ws.Range("L96").Value2 = 3.22
$ws.Range("M96").Value2 = "29/10/2023 13:59"
$ws.Range("N96").Value2 = 3.55
$ws.Range("O96").Value2 = "23/10/2023 19:12"
$ws.Range("P96").Value2 = 3.22
$ws.Range("Q96").Value2 = "29/10/2023 13:59"
$ws.Range("R96").Value2 = 2.71
$ws.Range("S96").Value2 = "23/10/2023 19:12"
$ws.Range("T96").Value2 = 2.46
$ws.Range("U96").Value2 = "29/10/2023 13:59"
$ws.Range("V96").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/f-karagumruk-trabzonspor/C03CnGT8/"

# Row 97
$ws.Range("A96:V96").Copy($ws.Range("A97"))
$ws.Range("B97").Value2 = "turkey"
$ws.Range("C97").Value2 = "super-lig"
$ws.Range("D97").Value2 = "2023-2024"
$ws.Range("A97").Value2 = 96
$ws.Range("E97").Value2 = 45228.70833333334
$ws.Range("F97").Value2 = "Ankaragucu"
$ws.Range("G97").Value2 = 2
$ws.Range("H97").Value2 = "Samsunspor"
$ws.Range("I97").Value2 = 0
$ws.Range("J97").Value2 = 2.19
$ws.Range("K97").Value2 = "23/10/2023 05:42"
$ws.Range("L97").Value2 = 2.7
$ws.Range("M97").Value2 = "29/10/2023 16:59"
$ws.Range("N97").Value2 = 3.59
$ws.Range("O97").Value2 = "23/10/2023 05:42"
$ws.Range("P97").Value2 = 3.32
$ws.Range("Q97").Value2 = "29/10/2023 16:54"
$ws.Range("R97").Value2 = 3.29
$ws.Range("S97").Value2 = "23/10/2023 05:42"
$ws.Range("T97").Value2 = 2.82
$ws.Range("U97").Value2 = "29/10/2023 16:59"
$ws.Range("V97").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-samsunspor/2kdPqEDR/"

# Row 98
$ws.Range("A97:V97").Copy($ws.Range("A98"))
$ws.Range("B98").Value2 = "turkey"
$ws.Range("C98").Value2 = "super-lig"
$ws.Range("D98").Value2 = "2023-2024"
$ws.Range("A98").Value2 = 97
$ws.Range("E98").Value2 = 45228.70833333334
$ws.Range("F98").Value2 = "Antalyaspor"
$ws.Range("G98").Value2 = 1
$ws.Range("H98").Value2 = "Basaksehir"
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 2.03
$ws.Range("K98").Value2 = "22/10/2023 20:15"
$ws.Range("L98").Value2 = 2.11
$ws.Range("M98").Value2 = "29/10/2023 16:54"
$ws.Range("N98").Value2 = 3.51
$ws.Range("O98").Value2 = "22/10/2023 20:15"
$ws.Range("P98").Value2 = 3.35
$ws.Range("Q98").Value2 = "29/10/2023 16:54"
$ws.Range("R98").Value2 = 3.8
$ws.Range("S98").Value2 = "22/10/2023 20:15"
$ws.Range("T98").Value2 = 3.89
$ws.Range("U98").Value2 = "29/10/2023 16:54"
$ws.Range("V98").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-basaksehir/tz8GozqF/"

# Row 99
$ws.Range("A98:V98").Copy($ws.Range("A99"))
$ws.Range("B99").Value2 = "turkey"
$ws.Range("C99").Value2 = "super-lig"
$ws.Range("D99").Value2 = "2023-2024"
$ws.Range("A99").Value2 = 98
$ws.Range("E99").Value2 = 45228.70833333334
$ws.Range("F99").Value2 = "Pendikspor"
$ws.Range("G99").Value2 = 0
$ws.Range("H99").Value2 = "Fenerbahce"
$ws.Range("I99").Value2 = 5
$ws.Range("J99").Value2 = 7.87
$ws.Range("K99").Value2 = "23/10/2023 05:42"
$ws.Range("L99").Value2 = 10.56
$ws.Range("M99").Value2 = "29/10/2023 16:59"
$ws.Range("N99").Value2 = 5.68
$ws.Range("O99").Value2 = "23/10/2023 05:42"
$ws.Range("P99").Value2 = 6.22
$ws.Range("Q99").Value2 = "29/10/2023 16:59"
$ws.Range("R99").Value2 = 1.35
$ws.Range("S99").Value2 = "23/10/2023 05:42"
$ws.Range("T99").Value2 = 1.28
$ws.Range("U99").Value2 = "29/10/2023 16:59"
$ws.Range("V99").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/pendikspor-fenerbahce/vc8IQY6k/"

# Row 100
$ws.Range("A99:V99").Copy($ws.Range("A100"))
$ws.Range("B100").Value2 = "turkey"
$ws.Range("C100").Value2 = "super-lig"
$ws.Range("D100").Value2 = "2023-2024"
$ws.Range("A100").Value2 = 99
$ws.Range("E100").Value2 = 45229.75
$ws.Range("F100").Value2 = "Besiktas"
$ws.Range("G100").Value2 = 2
$ws.Range("H100").Value2 = "Gaziantep"
$ws.Range("I100").Value2 = 0
$ws.Range("J100").Value2 = 1.34
$ws.Range("K100").Value2 = "23/10/2023 05:42"
$ws.Range("L100").Value2 = 1.56
$ws.Range("M100").Value2 = "30/10/2023 17:58"
$ws.Range("N100").Value2 = 5.57
$ws.Range("O100").Value2 = "23/10/2023 05:42"
$ws.Range("P100").Value2 = 4.38
$ws.Range("Q100").Value2 = "30/10/2023 17:58"
$ws.Range("R100").Value2 = 8.37
$ws.Range("S100").Value2 = "23/10/2023 05:42"
$ws.Range("T100").Value2 = 6.25
$ws.Range("U100").Value2 = "30/10/2023 17:58"
$ws.Range("V100").Value2 = "https://www.betexplorer.com/football/turkey/super-lig/besiktas-gaziantep/KQ6KpfbL/"

